# Regenerate s_vals data to filter save games.
# Updates TB (B), d2S (C), K (D), IP (E) and sum (G) columns
# for rows 2-6 on Sheet1 with the newly computed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2 = @{ B = 3.272327238179451; C = 1.626987699542094; D = 0.1496068669990043; E = 0.5333859586016987; G = 5.582307763322248 }
    3 = @{ B = 3.272327238179451; C = 1.626987699542094; D = 0.7210945179870265; E = 0.5333859586016987; G = 6.15379541431027 }
    4 = @{ B = 3.272327238179451; C = 1.626987699542094; D = 0.7210945179870265; E = 0.5333859586016987; G = 6.15379541431027 }
    5 = @{ B = 0.6545652718822623; C = 1.626987699542094; D = 0.1496068669990043; E = 0.5333859586016987; G = 2.964545797025059 }
    6 = @{ B = 3.272327238179451; C = 1.626987699542094; D = 0.7210945179870265; E = 0.5333859586016987; G = 6.15379541431027 }
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
